# Revise the metadata files
#
# 1) transport_solution lookup sheet: insert a new "PBS" option at the top,
#    and collapse "PBS (1x)" + "Formalin (10%; NBF)" into a single
#    "NBF (Neutral Buffered Formalin)" option (keeping its old slot).
# 2) .metadata sheet: bump pav:createdOn to the new timestamp.

$wb = $excel.ActiveWorkbook

$ts = $wb.Worksheets.Item("transport_solution")

$labels = @(
    "PBS",
    "Saline (Buffered)",
    "UWS",
    "DMEM",
    "Miltenyi Tissue Preservation Buffer",
    "NBF (Neutral Buffered Formalin)",
    "Unknown",
    "RPMI",
    "None",
    "HTK",
    "Belzer MPS/KPS"
)

$uris = @(
    "http://purl.obolibrary.org/obo/OBI_0100046",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000154",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000151",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185409",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000150",
    "http://purl.obolibrary.org/obo/OBIB_0000213",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C178973",
    "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000152",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000153"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ts.Cells.Item($row, 1).Value = $labels[$i]
    $ts.Cells.Item($row, 2).Value = $uris[$i]
}

$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2023-08-04T07:36:03-07:00"
